$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Homescreen" translation keys/rows below the existing table
# (rows 81-92), replicating the exact order in which the cells were
# authored so that shared-string indices come out the same.
$ws.Cells.Item(81, 1).Value = 'home_banner_title'
$ws.Cells.Item(81, 4).Value = 'మీ కలల ఆస్తిని కనుగొనండి'
$ws.Cells.Item(82, 1).Value = 'home_banner_subtitle'
$ws.Cells.Item(82, 2).Value = 'Explore premium real estate options in Vizag'
$ws.Cells.Item(82, 3).Value = 'विजाग में प्रीमियम रियल एस्टेट विकल्पों का अन्वेषण करें'
$ws.Cells.Item(82, 4).Value = 'విజాగ్‌లో ప్రీమియం రియల్ ఎస్టేట్ ఎంపికలను అన్వేషించండి'
$ws.Cells.Item(83, 1).Value = 'home_search_placeholder'
$ws.Cells.Item(83, 2).Value = 'Search properties...'
$ws.Cells.Item(83, 3).Value = 'संपत्तियाँ खोजें...'
$ws.Cells.Item(83, 4).Value = 'ఆస్తులను వెతకండి...'
$ws.Cells.Item(81, 2).Value = 'Find Your Dream Property'
$ws.Cells.Item(81, 3).Value = 'अपनी सपनों की संपत्ति खोजें'
$ws.Cells.Item(84, 1).Value = 'home_category_sites'
$ws.Cells.Item(84, 2).Value = 'Sites'
$ws.Cells.Item(84, 3).Value = 'साइट्स'
$ws.Cells.Item(84, 4).Value = 'సైట్లు'
$ws.Cells.Item(85, 1).Value = 'home_category_sites_desc'
$ws.Cells.Item(85, 2).Value = 'Plot Land'
$ws.Cells.Item(85, 3).Value = 'प्लॉट भूमि'
$ws.Cells.Item(85, 4).Value = 'ప్లాట్ భూమి'
$ws.Cells.Item(86, 1).Value = 'home_category_resorts'
$ws.Cells.Item(86, 2).Value = 'Resorts'
$ws.Cells.Item(86, 3).Value = 'रिसॉर्ट्स'
$ws.Cells.Item(86, 4).Value = 'రిసార్ట్స్'
$ws.Cells.Item(87, 1).Value = 'home_category_resorts_desc'
$ws.Cells.Item(87, 2).Value = 'Luxury Getaways'
$ws.Cells.Item(87, 3).Value = 'लक्ज़री अवकाश स्थल'
$ws.Cells.Item(87, 4).Value = 'లగ్జరీ విహార ప్రదేశాలు'
$ws.Cells.Item(88, 1).Value = 'home_category_flats'
$ws.Cells.Item(88, 2).Value = 'Flats'
$ws.Cells.Item(88, 3).Value = 'फ्लैट्स'
$ws.Cells.Item(88, 4).Value = 'ఫ్లాట్స్'
$ws.Cells.Item(89, 1).Value = 'home_category_flats_desc'
$ws.Cells.Item(89, 2).Value = 'Investment Land'
$ws.Cells.Item(89, 3).Value = 'निवेश भूमि'
$ws.Cells.Item(89, 4).Value = 'పెట్టుబడి భూమి'
$ws.Cells.Item(90, 1).Value = 'home_category_commercial'
$ws.Cells.Item(90, 2).Value = 'Commercial'
$ws.Cells.Item(90, 3).Value = 'व्यावसायिक'
$ws.Cells.Item(90, 4).Value = 'వాణిజ్య'
$ws.Cells.Item(91, 1).Value = 'home_category_commercial_desc'
$ws.Cells.Item(91, 2).Value = 'Business Spaces'
$ws.Cells.Item(91, 3).Value = 'व्यावसायिक स्थान'
$ws.Cells.Item(91, 4).Value = 'వ్యాపార స్థలాలు'
$ws.Cells.Item(92, 1).Value = 'home_language_modal_title'
$ws.Cells.Item(92, 2).Value = 'Select Your Language'
$ws.Cells.Item(92, 3).Value = 'अपनी भाषा चुनें'
$ws.Cells.Item(92, 4).Value = 'మీ భాషను ఎంచుకోండి'

# Update the view state (scroll position, zoom, active selection) to
# match the author's final view of the sheet.
$ws.Application.ActiveWindow.ScrollRow = 72
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.Zoom = 88
$ws.Range("B99").Select()
